$d = $word.ActiveDocument

# --- Paragraph 1: the placeholder/bookmark line -------------------------
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right) with 5pt "space" only
# (no line style/weight/color), matching <w:pBdr><w:top w:space="5"/>...).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5

# Increase the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25

# Replace the placeholder ID text and drop the trailing space run, so the
# paragraph ends up as a single run reading the new placeholder ID.
$p1 = $d.Paragraphs(1)
$textRange = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$textRange.Text = "**ID__AFFARS_AF_PGI_5342_202__ID**"
